# Append 5 new daily rows (29-10-2021 .. 02-11-2021) to the existing table,
# following the same pattern as the prior rows: column A = date label (text),
# column B = 3068, column C = 204.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("29-10-2021", "30-10-2021", "31-10-2021", "01-11-2021", "02-11-2021")

$startRow = 303
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)

    # Some of these strings (e.g. "01-11-2021") are ambiguous and would
    # otherwise be auto-converted to a date serial by Excel. Force them to
    # be stored as plain text, matching the rest of the column, then reset
    # the cell style so no extra formatting is left behind.
    $cellA.NumberFormat = "@"
    $cellA.Value = $dates[$i]
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = 3068
    $ws.Cells.Item($row, 3).Value = 204
}
